# Auto-generated cell updates applying the diff to the cryptos worksheet.
#
# Every cell in columns B-E on this sheet is a text cell in the source file
# (t="inlineStr"), even though some Price values look like plain numbers
# (e.g. "537.82", "0.0000139"). Assigning such a numeric-looking string to
# Range.Value makes Excel auto-convert it into a real number (changing the
# cell's type). To keep those specific cells as text we prepend a leading
# apostrophe - Excel's standard "force text" marker - which is stripped from
# the stored value while the cell keeps its original (General) number format.
# Values that are unambiguously text (names, links, percentages with a '%'
# sign) are assigned as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '58.921.57'; ForceText = $False }
    @{ Cell = 'D3'; Value = '2.529.20'; ForceText = $False }
    @{ Cell = 'E3'; Value = '  +3.39%  '; ForceText = $False }
    @{ Cell = 'E4'; Value = '  +0.03%  '; ForceText = $False }
    @{ Cell = 'D5'; Value = '537.82'; ForceText = $True }
    @{ Cell = 'E5'; Value = '  +0.29%  '; ForceText = $False }
    @{ Cell = 'D6'; Value = '144.06'; ForceText = $True }
    @{ Cell = 'E6'; Value = '  -2.43%  '; ForceText = $False }
    @{ Cell = 'D7'; Value = '0.998'; ForceText = $True }
    @{ Cell = 'E7'; Value = '  +0.12%  '; ForceText = $False }
    @{ Cell = 'D8'; Value = '0.571'; ForceText = $True }
    @{ Cell = 'E8'; Value = '  +0.24%  '; ForceText = $False }
    @{ Cell = 'D9'; Value = '2.525.18'; ForceText = $False }
    @{ Cell = 'E9'; Value = '  +2.69%  '; ForceText = $False }
    @{ Cell = 'D10'; Value = '0.0997'; ForceText = $True }
    @{ Cell = 'E10'; Value = '  +0.50%  '; ForceText = $False }
    @{ Cell = 'E11'; Value = '  +0.16%  '; ForceText = $False }
    @{ Cell = 'E12'; Value = '  +2.33%  '; ForceText = $False }
    @{ Cell = 'D13'; Value = '0.351'; ForceText = $True }
    @{ Cell = 'E13'; Value = '  +0.06%  '; ForceText = $False }
    @{ Cell = 'D14'; Value = '2.968.52'; ForceText = $False }
    @{ Cell = 'E14'; Value = '  +3.14%  '; ForceText = $False }
    @{ Cell = 'D15'; Value = '23.54'; ForceText = $True }
    @{ Cell = 'E15'; Value = '  -2.58%  '; ForceText = $False }
    @{ Cell = 'D16'; Value = '58.907.84'; ForceText = $False }
    @{ Cell = 'E16'; Value = '  -0.80%  '; ForceText = $False }
    @{ Cell = 'D17'; Value = '0.0000139'; ForceText = $True }
    @{ Cell = 'E17'; Value = '  +0.17%  '; ForceText = $False }
    @{ Cell = 'D18'; Value = '2.523.23'; ForceText = $False }
    @{ Cell = 'E18'; Value = '  +1.22%  '; ForceText = $False }
    @{ Cell = 'D19'; Value = '11.20'; ForceText = $True }
    @{ Cell = 'E19'; Value = '  +0.18%  '; ForceText = $False }
    @{ Cell = 'D20'; Value = '4.26'; ForceText = $True }
    @{ Cell = 'E20'; Value = '  -2.47%  '; ForceText = $False }
    @{ Cell = 'D21'; Value = '322.69'; ForceText = $True }
    @{ Cell = 'E21'; Value = '  -0.74%  '; ForceText = $False }
    @{ Cell = 'D22'; Value = '0.999'; ForceText = $True }
    @{ Cell = 'E22'; Value = '  +3.04%  '; ForceText = $False }
    @{ Cell = 'D23'; Value = '5.76'; ForceText = $True }
    @{ Cell = 'E23'; Value = '  +0.73%  '; ForceText = $False }
    @{ Cell = 'D24'; Value = '61.74'; ForceText = $True }
    @{ Cell = 'E24'; Value = '  +2.21%  '; ForceText = $False }
    @{ Cell = 'D25'; Value = '0.437'; ForceText = $True }
    @{ Cell = 'E25'; Value = '  -6.64%  '; ForceText = $False }
    @{ Cell = 'E26'; Value = '  +0.62%  '; ForceText = $False }
    @{ Cell = 'D27'; Value = '2.631.26'; ForceText = $False }
    @{ Cell = 'E27'; Value = '  +2.89%  '; ForceText = $False }
    @{ Cell = 'D28'; Value = '0.998'; ForceText = $True }
    @{ Cell = 'E28'; Value = '  +2.26%  '; ForceText = $False }
    @{ Cell = 'D29'; Value = '7.75'; ForceText = $True }
    @{ Cell = 'E29'; Value = '  +0.13%  '; ForceText = $False }
    @{ Cell = 'D30'; Value = '6.76'; ForceText = $True }
    @{ Cell = 'E30'; Value = '  -1.98%  '; ForceText = $False }
    @{ Cell = 'D31'; Value = '0.0₃0771'; ForceText = $False }
    @{ Cell = 'E31'; Value = '  +0.07%  '; ForceText = $False }
    @{ Cell = 'D32'; Value = '1.80'; ForceText = $True }
    @{ Cell = 'E32'; Value = '  -1.43%  '; ForceText = $False }
    @{ Cell = 'E33'; Value = '  -9.37%  '; ForceText = $False }
    @{ Cell = 'E34'; Value = '  +0.12%  '; ForceText = $False }
    @{ Cell = 'D35'; Value = '158.03'; ForceText = $True }
    @{ Cell = 'E35'; Value = '  +1.17%  '; ForceText = $False }
    @{ Cell = 'E36'; Value = '  +6.01%  '; ForceText = $False }
    @{ Cell = 'D37'; Value = '18.59'; ForceText = $True }
    @{ Cell = 'E37'; Value = '  +1.28%  '; ForceText = $False }
    @{ Cell = 'D38'; Value = '4.35'; ForceText = $True }
    @{ Cell = 'E38'; Value = '  -4.23%  '; ForceText = $False }
    @{ Cell = 'E39'; Value = '  -7.62%  '; ForceText = $False }
    @{ Cell = 'D40'; Value = '5.60'; ForceText = $True }
    @{ Cell = 'E40'; Value = '  -3.12%  '; ForceText = $False }
    @{ Cell = 'D41'; Value = '36.44'; ForceText = $True }
    @{ Cell = 'E41'; Value = '  -1.04%  '; ForceText = $False }
    @{ Cell = 'D42'; Value = '296.57'; ForceText = $True }
    @{ Cell = 'E42'; Value = '  -5.52%  '; ForceText = $False }
    @{ Cell = 'D43'; Value = '3.65'; ForceText = $True }
    @{ Cell = 'E43'; Value = '  -1.84%  '; ForceText = $False }
    @{ Cell = 'D44'; Value = '0.813'; ForceText = $True }
    @{ Cell = 'E44'; Value = '  -4.81%  '; ForceText = $False }
    @{ Cell = 'D45'; Value = '0.997'; ForceText = $True }
    @{ Cell = 'E45'; Value = '  +0.26%  '; ForceText = $False }
    @{ Cell = 'E46'; Value = '  +3.80%  '; ForceText = $False }
    @{ Cell = 'D47'; Value = '10.76'; ForceText = $True }
    @{ Cell = 'E47'; Value = '  +0.58%  '; ForceText = $False }
    @{ Cell = 'D48'; Value = '124.45'; ForceText = $True }
    @{ Cell = 'E48'; Value = '  +3.90%  '; ForceText = $False }
    @{ Cell = 'D49'; Value = '0.0930'; ForceText = $True }
    @{ Cell = 'E49'; Value = '  -1.03%  '; ForceText = $False }
    @{ Cell = 'E50'; Value = '  +0.53%  '; ForceText = $False }
    @{ Cell = 'B51'; Value = 'VeChain'; ForceText = $False }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $False }
    @{ Cell = 'D51'; Value = '0.0227'; ForceText = $True }
    @{ Cell = 'E51'; Value = '  -1.06%  '; ForceText = $False }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).Value = "'" + $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

